# Add files via upload
# Fills in sheet "2_" with a new multiple-choice-style question about
# differential equation solutions, updates the active sheet/tab selection
# state, and tweaks the selection left behind on the "6_MultAns" sheet.

$wb = $excel.ActiveWorkbook

# --- 6_MultAns : selection becomes the whole data range A1:C6 (instead
#           of the stray B30 cell). Do this before activating "2_" so the
#           final active tab stays "2_".
$wsAns = $wb.Worksheets.Item("6_MultAns")
$wsAns.Range("A1:C6").Select() | Out-Null

# --- 2_  : populate with the new question, answers and comments.
$ws = $wb.Worksheets.Item("2_")

# Column A (questions / prompts), column C (comments) -- written in this
# order so new shared-string entries land in the same order as the source
# workbook.
$ws.Range("A1").Value = "Imagine a complex differential equation that you don't know how to solve, like dx/dt = sin(e^t) + cos(x).  Which of the following is the *kind* of thing that could be a solution to the equation?"
$ws.Range("A2").Value = "A curve on a plot where every 't' value has only 1 'x' value"
$ws.Range("A4").Value = "x = 12t^2 + sin(6t)"
$ws.Range("A5").Value = "x = cos(At) + B sin(t^2) "

$ws.Range("C5").Value = "Yep!  This could be a general solution.   The function would be a solution no matter what A and B were."
$ws.Range("C4").Value = "Yep! This could be a particular solution: a single function that yields a true statement when it is plugged into the differential equation"
$ws.Range("C3").Value = "The solution to a differential equation is a function, not just a value."
$ws.Range("C2").Value = "Yep! This is just a visual representation of a function: the slope of this curve would be equal to sin(e^t) + cos(x) at every point on the curve."

# A3 is a lone numeric value (not text), matching the source data.
$ws.Range("A3").Value = 171.5

# Column B (Correct Y/N flags) re-uses existing shared strings.
$ws.Range("B1").Value = "Correct"
$ws.Range("B2").Value = "Y"
$ws.Range("B3").Value = "N"
$ws.Range("B4").Value = "Y"
$ws.Range("B5").Value = "Y"
$ws.Range("C1").Value = "Comment"

$ws.Rows.Item(1).RowHeight = 120
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 60

$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
